# Sample Project2 / Main.xlsx - SAVE
#
# The only semantically meaningful change in this commit is the value of
# cell C8 on the "Rules" sheet, which is bumped from 111 to 1111 (the rest
# of the captured diff is just column metadata noise re-emitted verbatim by
# the authoring tool on every save and carries no real content change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C8").Value = 1111.0
